$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.758.07"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.704.54"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "677.32"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.13"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.148"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.443"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000235"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.75"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.695.75"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.741.55"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.10"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.51"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "473.11"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.80"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.653"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.50"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.849.97"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000127"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.96"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.95"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.165"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.693.97"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0905"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.945"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.72"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.00"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.40"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000280"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.90"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "367.82"
$ws.Range("E51").Value = "  -1.18%  "
